{"js": "// Find the paragraph that starts with \"Prin aceast\u0103 activitate\" (the\n// introductory paragraph describing the activity) and collapse it down to\n// a single plain run containing the fully merged text, with the\n// constellation name updated from \"constela\u021bia Perseu\" to\n// \"Constela\u021bia Pegasus\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst searchKey = \"Prin aceast\u0103 activitate\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  if (text.indexOf(searchKey) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the target paragraph.\");\n}\n\nconst newText =\n  \"Prin aceast\u0103 activitate participa\u021bi \u00een cadrul unei campanii globale de \" +\n  \"observare \u0219i consemnare a celor mai slabe stele vizibile ca metod\u0103 de \" +\n  \"m\u0103surare a polu\u0103rii luminoase dintr-un anumit loc. Localiz\u00e2nd \u0219i \" +\n  \"observ\u00e2nd  Constela\u021bia Pegasus pe cerul nop\u021bii \u0219i compar\u00e2nd-o cu \" +\n  \"diagramele stelare, oamenii din \u00eentreaga lume vor putea afla \u00een ce \" +\n  \"m\u0103sur\u0103 iluminatul nocturn din comunitatea lor contribuie la poluarea \" +\n  \"luminoas\u0103. Contribu\u021biile dumneavoastr\u0103 la baza de date online vor \" +\n  \"facilita o documentare global\u0103 privind cerul nocturn observabil.\";\n\n// Clear out all the existing (multi-run) content/formatting of the\n// paragraph first so the replacement lands as a single run with no\n// leftover run-level formatting, then insert the merged text.\ntarget.clear();\ntarget.insertText(newText, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Collapse the introductory \"Prin aceast\u0103 activitate ...\" paragraph (which\n# is split across many small runs, one of them naming the constellation\n# \"Perseu\") down to a single plain run holding the fully merged text, with\n# the constellation name updated to \"Constela\u021bia Pegasus\".\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Prin aceast\u0103 activitate\")\nif (-not $found) {\n    throw \"Could not locate the target paragraph.\"\n}\n\n# Expand the found hit out to its whole paragraph, then pull the end back\n# one character so the paragraph mark itself is not touched (deleting it\n# would merge this paragraph with the next one).\n$rng.Expand(4)        # wdParagraph\n$rng.MoveEnd(1, -1)   # wdCharacter, pull back off the paragraph mark\n$rng.Delete()\n\n$newText = \"Prin aceast\u0103 activitate participa\u021bi \u00een cadrul unei campanii globale de observare \u0219i consemnare a celor mai slabe stele vizibile ca metod\u0103 de m\u0103surare a polu\u0103rii luminoase dintr-un anumit loc. Localiz\u00e2nd \u0219i observ\u00e2nd  Constela\u021bia Pegasus pe cerul nop\u021bii \u0219i compar\u00e2nd-o cu diagramele stelare, oamenii din \u00eentreaga lume vor putea afla \u00een ce m\u0103sur\u0103 iluminatul nocturn din comunitatea lor contribuie la poluarea luminoas\u0103. Contribu\u021biile dumneavoastr\u0103 la baza de date online vor facilita o documentare global\u0103 privind cerul nocturn observabil.\"\n\n$rng.InsertBefore($newText)\n"}
